$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Paragraph 24: "{%p if (users[0].insurance | length > 5) or
#    (other_parties[0].insurance | length > 5) %}" becomes
#    "{%p if there_are_marital_children and ((users[0].insurance |
#    length > 5) or (other_parties[0].insurance | length > 5)) %}"
# ------------------------------------------------------------------
$p = $d.Paragraphs.Item(24)
$r = $p.Range
$r.Find.Execute("{%p if (", $true, $false, $false, $false, $false, $true, 1, $false, `
    "{%p if there_are_marital_children and ((", 2)

# Re-fetch the paragraph range (text length changed) and fix the tail.
$p = $d.Paragraphs.Item(24)
$r = $p.Range
$r.Find.Execute("length > 5) %}", $true, $false, $false, $false, $false, $true, 1, $false, `
    "length > 5)) %}", 2)

# ------------------------------------------------------------------
# 2. Remove the "{%p endif %}" paragraph that used to immediately
#    follow "MEDICAL SUPPORT PROVISIONS, continued" (paragraph 26).
# ------------------------------------------------------------------
$endifPara = $d.Paragraphs.Item(26)
$endifPara.Range.Delete()

# ------------------------------------------------------------------
# 3. At the end of the document, insert a new "{%p endif %}" paragraph
#    (matching the one removed above, keepNext + same formatting)
#    right before the final "{%p endif %}" paragraph, and append a new
#    blank paragraph right after that final paragraph.
# ------------------------------------------------------------------
$lastCount = $d.Paragraphs.Count
$finalEndif = $d.Paragraphs.Item($lastCount)

# Insert the new endif paragraph before the final one.
$finalEndif.Range.InsertParagraphBefore()

# The newly inserted (now empty) paragraph is the one that used to be
# at index $lastCount; fill it in with the endif text and keepNext.
$newPara = $d.Paragraphs.Item($lastCount)
$newPara.Range.Text = "{%p endif %}"
$newPara.Format.KeepWithNext = $true

# Re-fetch the final paragraph (index shifted by +1) and append a
# trailing blank paragraph after it.
$finalEndif = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalEndif.Range.InsertParagraphAfter()
